$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 24.60000000000041
$ws.Range("G2").Value = 0.003706323675013001
$ws.Range("H2").Value = 0.009905508549801585
$ws.Range("K2").Value = 4.809352232206899
$ws.Range("L2").Value = "[1.5712018191371477, 8.04750264527665]"
$ws.Range("M2").Value = 0.003768162521081253
$ws.Range("N2").Value = 0.003768162521081253
$ws.Range("O2").Value = -1.559789745926464
$ws.Range("P2").Value = "[-2.4906320136567732, -0.6289474781961548]"
$ws.Range("Q2").Value = 0.001111521131287763
$ws.Range("R2").Value = 0.001111521131287763
$ws.Range("S2").Value = 13.59673130643606
$ws.Range("T2").Value = "[11.604952510276265, 15.588510102595853]"
$ws.Range("W2").Value = 6.106906906907007
$ws.Range("X2").Value = 2.462462462462502
$ws.Range("Y2").Value = 9.751351351351513
$ws.Range("B3").Value = 1
$ws.Range("E3").Value = 24.99000000000047
$ws.Range("G3").Value = 0.0004315006780571951
$ws.Range("H3").Value = 0.004232340223790355
$ws.Range("K3").Value = 5.339318214720941
$ws.Range("L3").Value = "[1.8052935063489617, 8.87334292309292]"
$ws.Range("M3").Value = 0.00319284221458771
$ws.Range("N3").Value = 0.003768162521081253
$ws.Range("O3").Value = -2.717053105807389
$ws.Range("P3").Value = "[-3.421474281387082, -2.012631930227696]"
$ws.Range("Q3").Value = [double]"4.518607710224387e-13"
$ws.Range("R3").Value = [double]"9.037215420448774e-13"
$ws.Range("S3").Value = 13.38596429365885
$ws.Range("T3").Value = "[11.497878553023913, 15.27405003429379]"
$ws.Range("W3").Value = 10.80648648648669
$ws.Range("X3").Value = 8.004804804804959
$ws.Range("Y3").Value = 13.60816816816842
